# Update xgboost_params sheet: re-order/rename the parameter-set headers,
# collapse from 8 data columns (B:I) down to 5 (B:F), and update the
# numeric values to match the new column layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 headers (B1:F1) ---
$ws.Range("B1").Value = "H2n"
$ws.Range("C1").Value = "iron VDZP"
$ws.Range("D1").Value = "ozone MB"
$ws.Range("E1").Value = "ozone VDZP"
$ws.Range("F1").Value = "ozone VTZP"

# --- Row 2: objective (unchanged values, just fewer columns now) ---
$ws.Range("B2").Value = "reg:squarederror"
$ws.Range("C2").Value = "reg:squarederror"
$ws.Range("D2").Value = "reg:squarederror"
$ws.Range("E2").Value = "reg:squarederror"
$ws.Range("F2").Value = "reg:squarederror"

# --- Row 3: enable_categorical ---
$ws.Range("B3").Value = $false
$ws.Range("C3").Value = $false
$ws.Range("D3").Value = $false
$ws.Range("E3").Value = $false
$ws.Range("F3").Value = $false

# --- Row 4: max_depth ---
$ws.Range("B4").Value = 10
$ws.Range("C4").Value = 100
$ws.Range("D4").Value = 10
$ws.Range("E4").Value = 100
$ws.Range("F4").Value = 10

# --- Row 5: n_estimators ---
$ws.Range("B5").Value = 100
$ws.Range("C5").Value = 100
$ws.Range("D5").Value = 100
$ws.Range("E5").Value = 100
$ws.Range("F5").Value = 100

# --- Row 6: reg_alpha ---
$ws.Range("B6").Value = 0.000001
$ws.Range("C6").Value = 0.001
$ws.Range("D6").Value = 0.001
$ws.Range("E6").Value = 0.1
$ws.Range("F6").Value = 0.001

# --- Row 7: reg_lambda ---
$ws.Range("B7").Value = 0.001
$ws.Range("C7").Value = 0.1
$ws.Range("D7").Value = 0.1
$ws.Range("E7").Value = 0.001
$ws.Range("F7").Value = 0.001

# --- Remove the now-unused columns G:I entirely (data + formatting) ---
$ws.Range("G1:I7").Delete()
